# Refresh "想去人数" (column F) and "最低票价" (column G) figures
# to the latest scraped values (gh-pages data output).
$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$wsExpo = $wb.Worksheets.Item("展览")
$data = @(
    @(2, 1158, -1),
    @(3, 882, 60),
    @(5, 54, -1),
    @(6, 1115, -1),
    @(8, 2397, -1),
    @(9, 7826, -1),
    @(10, 932, -1),
    @(11, 454, -1),
    @(12, 395, -1),
    @(14, 435, -1),
    @(16, 165, -1),
    @(17, 8040, -1),
    @(19, 1396, -1),
    @(20, 160, -1),
    @(25, 178, -1),
    @(28, 114, -1),
    @(29, 32, -1),
    @(31, 1164, -1),
    @(33, 101, -1),
    @(34, 68, -1),
    @(35, 87, -1),
    @(36, 46, -1),
    @(37, 82, -1)
)

foreach ($item in $data) {
    $r = $item[0]
    $f = $item[1]
    $g = $item[2]
    $wsExpo.Cells.Item($r, 6).Value = $f
    if ($g -ge 0) {
        $wsExpo.Cells.Item($r, 7).Value = $g
    }
}

# --- Sheet "全部类型" (All types) mirrors the same rows ---
$wsAll = $wb.Worksheets.Item("全部类型")
$data = @(
    @(2, 1158, -1),
    @(3, 882, 60),
    @(5, 54, -1),
    @(6, 1115, -1),
    @(8, 2397, -1),
    @(9, 7826, -1),
    @(10, 932, -1),
    @(11, 454, -1),
    @(12, 395, -1),
    @(14, 435, -1),
    @(16, 165, -1),
    @(17, 8040, -1),
    @(19, 1396, -1),
    @(20, 160, -1),
    @(25, 178, -1),
    @(28, 114, -1),
    @(29, 32, -1),
    @(31, 1164, -1),
    @(34, 68, -1),
    @(35, 87, -1),
    @(36, 46, -1),
    @(37, 82, -1)
)

foreach ($item in $data) {
    $r = $item[0]
    $f = $item[1]
    $g = $item[2]
    $wsAll.Cells.Item($r, 6).Value = $f
    if ($g -ge 0) {
        $wsAll.Cells.Item($r, 7).Value = $g
    }
}
